# Amend corrected label annotations
# Normalize the "labels" column (F) text casing / ordering per the corrected
# label annotations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    17  = "18_hazards_to_humans_and_domestic_animals"
    19  = "ppe"
    24  = "env warning - water || off target movement"
    34  = "application instructions"
    35  = "application instructions"
    36  = "135_product_information"
    71  = "use restrictions"
    72  = "use restrictions"
    73  = "use restrictions || irrigation"
    74  = "use restrictions"
    75  = "off target movement"
    76  = "off target movement"
    77  = "off target movement"
    78  = "172_sensitive_areas"
    80  = "172_sensitive_areas || off target movement"
    82  = "application instructions"
    83  = "mixing || application instructions"
    84  = "application instructions"
    85  = "mixing || application instructions"
    86  = "application instructions"
    88  = "mixing"
    89  = "mixing"
    90  = "mixing"
    91  = "mixing"
    92  = "mixing"
    93  = "mixing"
    94  = "mixing"
    95  = "mixing"
    96  = "mixing"
    98  = "safety procedures"
    99  = "mixing"
    119 = "mixing"
    121 = "mixing"
    123 = "154_pesticide_storage"
    135 = "18_hazards_to_humans_and_domestic_animals"
    145 = "env warning - water || off target movement"
    146 = "154_pesticide_storage"
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
